# Unificación de logs, mains y modularización
# Add two new data rows (8 and 9) to the existing "servicios_internet" sheet,
# replicating the "Personal" / "300" offer with two different price text
# representations ("24000" and "24.000").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Personal"
$ws.Range("B8").Value = "'300"
$ws.Range("C8").Value = "'24000"
$ws.Range("A8:C8").Style = "Normal"

$ws.Range("A9").Value = "Personal"
$ws.Range("B9").Value = "'300"
$ws.Range("C9").Value = "'24.000"
$ws.Range("A9:C9").Style = "Normal"
